$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Replace the bold "REF" placeholder (before "based on how
#    representative...") with the Minderman/Park citation text, with
#    "et al." in italics.
# ---------------------------------------------------------------------
$refRng = $d.Content
$foundRef = $refRng.Find.Execute("REF", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundRef) {
    $refRng.Font.Bold = $false
    $base = $refRng.Start

    $part1 = "(Minderman"
    $part2 = " "
    $part3 = "et al."
    $part4 = " "
    $part5 = "2012; Park, Turner & Minderman 2013)"

    $refRng.Text = $part1 + $part2 + $part3 + $part4 + $part5

    $iStart = $base + $part1.Length + $part2.Length
    $iEnd = $iStart + $part3.Length
    $italicRng = $d.Range($iStart, $iEnd)
    $italicRng.Font.Italic = $true
}

# ---------------------------------------------------------------------
# 2. "Bat data" -> "Bat data and transects" (heading + bookmark rename)
# ---------------------------------------------------------------------
$batBm = $d.Bookmarks("bat-data")
$d.Bookmarks.Add("bat-data-and-transects", $batBm.Range)
$d.Bookmarks("bat-data").Delete()

$batRng = $d.Content
$foundBat = $batRng.Find.Execute("Bat data", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundBat) {
    $batRng.Text = "Bat data and transects"
}

# ---------------------------------------------------------------------
# 3. "Measure of bat activity: probability of a pass per hectare
#    surveyed" -> "Bat activity: probability of a pass per hectare
#    surveyed" (heading + bookmark rename)
# ---------------------------------------------------------------------
$measureBm = $d.Bookmarks("measure-of-bat-activity-probability-of-a-pass-per-hectare-surveyed")
$d.Bookmarks.Add("bat-activity-probability-of-a-pass-per-hectare-surveyed", $measureBm.Range)
$d.Bookmarks("measure-of-bat-activity-probability-of-a-pass-per-hectare-surveyed").Delete()

$measureRng = $d.Content
$foundMeasure = $measureRng.Find.Execute("Measure of bat activity: probability of a pass per hectare surveyed", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundMeasure) {
    $measureRng.Text = "Bat activity: probability of a pass per hectare surveyed"
}

# ---------------------------------------------------------------------
# 4. Reword the "skewed counts" sentence.
# ---------------------------------------------------------------------
$skewOld = "was so skewed (many zeros and excessive variation) that count-based statistical models did not provide any reasonable fit. Second, using bat 'passes' as a measure of activity provides only a relative measure of activity in he first place."
$skewNew = "was highly skewed (many zeros and excessive variation) so that count-based statistical models did not provide any reasonable fit. Second, using bat 'passes' as a measure of activity provides a relative measure of activity and analyses of absolute pass count would therefore add little information."

$skewRng = $d.Content
$foundSkew = $skewRng.Find.Execute($skewOld, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundSkew) {
    $skewRng.Text = $skewNew
}

# ---------------------------------------------------------------------
# 5. "References" heading: Heading2 -> Heading1
# ---------------------------------------------------------------------
$refsRng = $d.Content
$foundRefs = $refsRng.Find.Execute("References", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundRefs) {
    $refsPara = $refsRng.Paragraphs.First
    $refsPara.Style = "Heading 1"
}

# ---------------------------------------------------------------------
# 6. Add "Appendix 1" heading (Heading1) with bookmark, after References
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$appendixPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$appendixPara.Style = "Heading 1"
$appendixPara.Range.Text = "Appendix 1"
$appendixBmStart = $d.Range($appendixPara.Range.Start, $appendixPara.Range.Start)
$d.Bookmarks.Add("appendix-1", $appendixBmStart)

# ---------------------------------------------------------------------
# 7. Add bibliography entry: Minderman et al. (2012)
# ---------------------------------------------------------------------
$appendixPara.Range.InsertParagraphAfter()
$bib1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$bib1.Style = "Bibliography"

$b1p1 = "Minderman, J., Pendlebury, C.J., Pearce-Higgins, J.W. & Park, K.J. (2012) Experimental Evidence for the Effect of Small Wind Turbine Proximity and Operation on Bird and Bat Activity."
$b1p2 = " "
$b1p3 = "PLoS ONE"
$b1p4 = ","
$b1p5 = " "
$b1p6 = "7"
$b1p7 = ", e41177."

$bib1.Range.Text = $b1p1 + $b1p2 + $b1p3 + $b1p4 + $b1p5 + $b1p6 + $b1p7
$bib1Base = $bib1.Range.Start

$b1IStart = $bib1Base + $b1p1.Length + $b1p2.Length
$b1IEnd = $b1IStart + $b1p3.Length
$d.Range($b1IStart, $b1IEnd).Font.Italic = $true

$b1BStart = $b1IEnd + $b1p4.Length + $b1p5.Length
$b1BEnd = $b1BStart + $b1p6.Length
$d.Range($b1BStart, $b1BEnd).Font.Bold = $true

# ---------------------------------------------------------------------
# 8. Add bibliography entry: Park, Turner & Minderman (2013)
# ---------------------------------------------------------------------
$bib1.Range.InsertParagraphAfter()
$bib2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$bib2.Style = "Bibliography"

$b2p1 = "Park, K.J., Turner, A. & Minderman, J. (2013) Integrating applied ecology and planning policy: The case of micro-turbines and wildlife conservation."
$b2p2 = " "
$b2p3 = "Journal of Applied Ecology"
$b2p4 = ","
$b2p5 = " "
$b2p6 = "50"
$b2p7 = ", 199–204."

$bib2.Range.Text = $b2p1 + $b2p2 + $b2p3 + $b2p4 + $b2p5 + $b2p6 + $b2p7
$bib2Base = $bib2.Range.Start

$b2IStart = $bib2Base + $b2p1.Length + $b2p2.Length
$b2IEnd = $b2IStart + $b2p3.Length
$d.Range($b2IStart, $b2IEnd).Font.Italic = $true

$b2BStart = $b2IEnd + $b2p4.Length + $b2p5.Length
$b2BEnd = $b2BStart + $b2p6.Length
$d.Range($b2BStart, $b2BEnd).Font.Bold = $true

Write-Output "Edits applied"
